# Text updates as supplied by PM&C.
# Split the combined "Source" note on the Description sheet into a
# labelled "Source" row plus two separate sentences (Indigenous /
# Non-Indigenous data sources), matching the formatting already used
# for the other footnote rows on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Description")

$indigenous = "Indigenous data sourced from the ABS National Aboriginal and Torres Strait Islander Social Survey (NATSISS) (2008 and 2014-15)."
$nonIndigenous = "Non-Indigenous data sourced from the ABS Survey of Education and Work (SEW) (2008 and 2014 )."

# Row 10 currently holds the combined note in B10 (A10 empty). Give it a
# "Source" label in A10 and keep the Indigenous sentence in B10, matching
# the text-wrap style (Arial 12, same dark colour) used by the other note
# rows (B5:B9).
$ws.Range("A10").Value = "Source"
$ws.Range("B10").Value = $indigenous
$ws.Range("B10").Font.Size = 12
$ws.Range("B10").Font.Color = 655360

# Insert the Non-Indigenous sentence as a new row 11 in column B, using
# the same wrapped style.
$ws.Range("B11").Value = $nonIndigenous
$ws.Range("B11").Font.Size = 12
$ws.Range("B11").Font.Color = 655360
